$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.609.82"
$ws.Range("E2").Value = "  +1.61%  "

# Row 3
$ws.Range("D3").Value = "1.803.93"
$ws.Range("E3").Value = "  +0.89%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'224.16"
$ws.Range("E5").Value = "  -1.61%  "

# Row 6
$ws.Range("E6").Value = "  -0.15%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("E8").Value = "  +3.47%  "

# Row 9
$ws.Range("E9").Value = "  +2.92%  "

# Row 10
$ws.Range("D10").Value = "'0.0712"
$ws.Range("E10").Value = "  +7.94%  "

# Row 11
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "  -0.01%  "

# Row 12
$ws.Range("E12").Value = "  +0.94%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.839.64"
$ws.Range("E13").Value = "  +2.98%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.09"
$ws.Range("E14").Value = "  -3.62%  "

# Row 15
$ws.Range("E15").Value = "  +1.07%  "

# Row 16
$ws.Range("D16").Value = "34.650.69"
$ws.Range("E16").Value = "  +1.75%  "

# Row 17
$ws.Range("E17").Value = "  +2.16%  "

# Row 18
$ws.Range("D18").Value = "'69.22"
$ws.Range("E18").Value = "  -0.46%  "

# Row 19
$ws.Range("D19").Value = "'252.79"
$ws.Range("E19").Value = "  -0.15%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0803"
$ws.Range("E20").Value = "  +8.13%  "

# Row 21
$ws.Range("D21").Value = "'11.09"
$ws.Range("E21").Value = "  +5.91%  "

# Row 23
$ws.Range("E23").Value = "  -0.42%  "

# Row 25
$ws.Range("D25").Value = "'161.45"
$ws.Range("E25").Value = "  +2.91%  "

# Row 26
$ws.Range("E26").Value = "  -0.96%  "

# Row 27
$ws.Range("D27").Value = "'7.15"
$ws.Range("E27").Value = "  +1.64%  "

# Row 28
$ws.Range("E28").Value = "  +0.00%  "

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0528"
$ws.Range("E30").Value = "  +2.26%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.80"
$ws.Range("E31").Value = "  -0.35%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.21"
$ws.Range("E32").Value = "  -0.26%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'3.62"
$ws.Range("E33").Value = "  +0.29%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.89"
$ws.Range("E34").Value = "  +2.55%  "

# Row 35
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.433.36"
$ws.Range("E35").Value = "  -1.37%  "

# Row 36
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.07"
$ws.Range("E36").Value = "  -0.10%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.644"
$ws.Range("E37").Value = "  +2.37%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0192"
$ws.Range("E38").Value = "  +2.63%  "

# Row 39
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").Value = "'84.83"
$ws.Range("E39").Value = "  +1.53%  "

# Row 40
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.960"
$ws.Range("E40").Value = "  +6.48%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.80"
$ws.Range("E41").Value = "  -0.85%  "

# Row 42
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.35"
$ws.Range("E42").Value = "  +0.10%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.16"
$ws.Range("E43").Value = "  +3.89%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.04"
$ws.Range("E44").Value = "  +4.44%  "

# Row 45
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.06"
$ws.Range("E45").Value = "  -1.04%  "

# Row 46
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").Value = "'0.0498"
$ws.Range("E46").Value = "  -2.77%  "

# Row 47
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.958.66"
$ws.Range("E47").Value = "  +0.57%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'12.32"
$ws.Range("E48").Value = "  +3.32%  "

# Row 49
$ws.Range("D49").Value = "'106.52"
$ws.Range("E49").Value = "  +8.69%  "

# Row 50
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.08%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0126"
$ws.Range("E51").Value = "  +8.68%  "
